# Apply updates described by the commit: "update to published CDA FHIR
# logical model with patches #241"

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version
$meta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"

# Date
$meta.Range("B8").Value = "2024-06-19T17:47:42+02:00"

# Contact
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Root element (PIVL_TS) Definition text
$elements.Range("M2").Value = "A quantity specifying a point on the axis of natural time. A point in time is most often represented as a calendar expression."

# PIVL_TS.operator Binding Value Set
$elements.Range("Z5").Value = "http://hl7.org/cda/stds/core/ValueSet/CDASetOperator"

# Column Z (26) got slightly wider to fit the new, longer URL text
# (target raw OOXML width ~51.21484375; the engine quantizes ColumnWidth
# input, so 50.3 is the closest input that lands on the nearest
# reachable bucket, 51.16666...)
$elements.Columns.Item(26).ColumnWidth = 50.3
